$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 11 entirely (a row not present upstream); rows below shift up.
$ws.Rows("11:11").Delete()
